# Fixed a tbd search bug
#
# Appends the list of orphan-tag paragraphs after the intro paragraph
# ("These are the orphan tags that were found in the documents: ").
# Each new paragraph is a single, unformatted run so the inserted
# markup matches plain `<w:r><w:t>...</w:t></w:r>` (no inherited bold
# from the preceding run). We build the paragraph OOXML by hand and
# insert it with Range.InsertXML so no run/paragraph formatting is
# carried over from the insertion point.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$orphanTags = @(
    @{ Text = "TARGEST:UI:300 ";    Preserve = $true  },
    @{ Text = "TARGEST:UI:500 ";    Preserve = $true  },
    @{ Text = "TARGEST:HLR:100 ";   Preserve = $true  },
    @{ Text = "TARGEST:HLR:200 ";   Preserve = $true  },
    @{ Text = "TARGEST:HLR:300 ";   Preserve = $true  },
    @{ Text = "TARGEST:TBD:100";    Preserve = $false },
    @{ Text = "TARGEST:SYS:500";    Preserve = $false },
    @{ Text = ".TARGEST:FUNC:200";  Preserve = $false }
)

foreach ($tag in $orphanTags) {
    if ($tag.Preserve) {
        $tOpen = '<w:t xml:space="preserve">'
    } else {
        $tOpen = '<w:t>'
    }

    $paraXml = "<w:p $wNs><w:r>$tOpen$($tag.Text)</w:t></w:r></w:p>"

    # Insert right before the very end of the document body (the end
    # of the final paragraph mark), so each new paragraph lands after
    # the previous one, right before the sectPr.
    $insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
    $insertionPoint.InsertXML($paraXml)
}
